$p = $ppt.ActivePresentation

# Slide 4: "Hub Design" -> "Timeline"
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Timeline"

# Slide 5: "Node Design" -> "Budget"
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Budget"

# Slide 6: "FAQs" -> "Questions?" and remove the "Doth we even need this" placeholder
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Questions?"
# Deleting a layout-backed placeholder once only resets it to its empty
# default state (PowerPoint regenerates it from the layout); delete again
# to actually remove the now-empty placeholder shape from the slide.
$s6.Shapes.Item(2).Delete()
$s6.Shapes.Item(2).Delete()
